$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D2:D51) to be treated as text so values like
# '1.002' or '0.7839' are not auto-converted to numbers by Excel, matching
# the source data which stores these as inline/shared strings without a
# custom number format.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# --- Column D (Price) updates ---
$ws.Range('D2').Value = '31.222.78'
$ws.Range('D3').Value = '1.994.54'
$ws.Range('D4').Value = '1.002'
$ws.Range('D5').Value = '0.7839'
$ws.Range('D6').Value = '254.50'
$ws.Range('D7').Value = '1.002'
$ws.Range('D8').Value = '0.3477'
$ws.Range('D9').Value = '27.74'
$ws.Range('D10').Value = '0.07016'
$ws.Range('D11').Value = '0.8432'
$ws.Range('D12').Value = '0.08179'
$ws.Range('D13').Value = '100.49'
$ws.Range('D14').Value = '1.991.42'
$ws.Range('D15').Value = '5.614'
$ws.Range('D16').Value = '15.37'
$ws.Range('D17').Value = '272.38'
$ws.Range('D18').Value = '31.238.32'
$ws.Range('D19').Value = '0.000008010'
$ws.Range('D20').Value = '5.864'
$ws.Range('D21').Value = '2.251.77'
$ws.Range('D22').Value = '1.001'
$ws.Range('D23').Value = '1.001'
$ws.Range('D24').Value = '7.046'
$ws.Range('D25').Value = '10.01'
$ws.Range('D26').Value = '166.28'
$ws.Range('D27').Value = '0.1462'
$ws.Range('D28').Value = '19.87'
$ws.Range('D29').Value = '2.339'
$ws.Range('D30').Value = '1.598'
$ws.Range('D31').Value = '1.361'
$ws.Range('D32').Value = '4.595'
$ws.Range('D33').Value = '4.429'
$ws.Range('D34').Value = '0.05260'
$ws.Range('D35').Value = '0.7814'
$ws.Range('D36').Value = '1.223'
$ws.Range('D37').Value = '2.757'
$ws.Range('D38').Value = '1.001'
$ws.Range('D39').Value = '0.02004'
$ws.Range('D40').Value = '2.904'
$ws.Range('D41').Value = '6.725'
$ws.Range('D42').Value = '79.52'
$ws.Range('D43').Value = '0.4667'
$ws.Range('D44').Value = '2.106'
$ws.Range('D45').Value = '0.8535'
$ws.Range('D46').Value = '104.71'
$ws.Range('D47').Value = '1.001'
$ws.Range('D48').Value = '10.05'
$ws.Range('D49').Value = '7.655'
$ws.Range('D50').Value = '1.593'
$ws.Range('D51').Value = '37.31'

# Remove the temporary text-format styling so the cells keep the same
# (unstyled) appearance as before, only the underlying text changed.
$priceRange.ClearFormats()

# --- Other column updates (Coin name, Link, Volume) ---
$ws.Range('E2').Value = '  +2.09%  '
$ws.Range('E3').Value = '  +5.90%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('E5').Value = '  +65.37%  '
$ws.Range('E6').Value = '  +3.30%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('E8').Value = '  +20.00%  '
$ws.Range('E9').Value = '  +23.72%  '
$ws.Range('E10').Value = '  +7.34%  '
$ws.Range('E11').Value = '  +9.58%  '
$ws.Range('E12').Value = '  +4.54%  '
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('E14').Value = '  +5.72%  '
$ws.Range('E15').Value = '  +7.07%  '
$ws.Range('E16').Value = '  +16.21%  '
$ws.Range('E17').Value = '  -4.45%  '
$ws.Range('E18').Value = '  +2.24%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('E19').Value = '  +6.42%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('E20').Value = '  +9.62%  '
$ws.Range('E21').Value = '  +5.64%  '
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('E24').Value = '  +9.32%  '
$ws.Range('E25').Value = '  +9.01%  '
$ws.Range('E26').Value = '  +1.69%  '
$ws.Range('E27').Value = '  +50.72%  '
$ws.Range('E28').Value = '  +4.02%  '
$ws.Range('E29').Value = '  +22.86%  '
$ws.Range('E30').Value = '  +6.49%  '
$ws.Range('E31').Value = '  +1.73%  '
$ws.Range('E32').Value = '  +8.02%  '
$ws.Range('E33').Value = '  +5.86%  '
$ws.Range('E34').Value = '  +8.57%  '
$ws.Range('E35').Value = '  +11.74%  '
$ws.Range('E36').Value = '  +8.31%  '
$ws.Range('E37').Value = '  +0.49%  '
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('E39').Value = '  +4.95%  '
$ws.Range('E40').Value = '  +0.77%  '
$ws.Range('E41').Value = '  +6.71%  '
$ws.Range('E42').Value = '  +4.79%  '
$ws.Range('E43').Value = '  +9.56%  '
$ws.Range('E44').Value = '  +6.24%  '
$ws.Range('E45').Value = '  +2.66%  '
$ws.Range('E46').Value = '  +3.25%  '
$ws.Range('E47').Value = '  +0.09%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E48').Value = '  +1.56%  '
$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('E49').Value = '  +8.88%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('E50').Value = '  +18.70%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('E51').Value = '  +6.12%  '

